# Update PLC data 2025-10-13 13:49:00
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 161154
$ws.Range("C4").Value = 152172
$ws.Range("C5").Value = 8982
$ws.Range("C8").Value = 64.53
